$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 3. Số ngày làm việc theo quy định của pháp luật lao động trong tháng: 0 -> 23
$ws.Range("A9").Value = "3. Số ngày làm việc theo quy định của pháp luật lao động trong tháng: 23"

# 7. Hành vi vi phạm:  ->  7. Hành vi vi phạm: 0
$ws.Range("F12").Value = "7. Hành vi vi phạm: 0"

# 8. Hình thức kỷ luật:  ->  8. Hình thức kỷ luật: 0
$ws.Range("I12").Value = "8. Hình thức kỷ luật: 0"

# Update the details of the first work item row (row 16)
$ws.Range("C16").Value = "26/03/2025"
$ws.Range("D16").Value = "26/03/2025"
$ws.Range("E16").Value = "26/03/2025"
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = "updating..."
$ws.Range("H16").Value = "Vượt tiến độ hoặc có chất lượng hoặc hiệu quả cao"
$ws.Range("I16").Value = "Đúng tiến độ, đảm bảo chất lượng và hiệu quả"
$ws.Range("J16").Value = "Nguyễn Văn Hải"
$ws.Range("K16").Value = "Chưa đảm bảo về yêu cầu và chất lượng"
$ws.Range("L16").Value = "Nguyễn Duy Cường"

# Remove the two other work item rows (old rows 17 and 18); this shifts the
# trailing "10. Kết quả..." / "Cán bộ lập phiếu" rows up from 20/21 to 18/19
# and shrinks the used range from A1:M21 to A1:M19 automatically.
$ws.Rows("17:18").Delete()

# Restore the selection to match the updated table extent
$ws.Range("A15:M16").Select()
